$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '408-7470107-3753964'
$ws.Range("C2").Value = '2024-11-22T16:08:05+00:00'
$ws.Range("D2").Value = '2024-11-22T16:22:54+00:00'
$ws.Range("F2").Value = 'KAR MARK kudampuli lehyam 1KG'
$ws.Range("I2").Value = 1230
$ws.Range("J2").Value = 131.79
$ws.Range("L2").Value = 8.57

$ws.Range("B3").Value = '407-3594126-7517156'
$ws.Range("C3").Value = '2024-11-22T15:28:15+00:00'
$ws.Range("D3").Value = '2024-11-22T18:18:58+00:00'
$ws.Range("F3").Value = 'KAR MARK kudampuli lehyam 1KG'
$ws.Range("I3").Value = 1230
$ws.Range("J3").Value = 131.78

$ws.Range("B4").Value = '406-9842547-5707502'
$ws.Range("C4").Value = '2024-11-22T14:02:57+00:00'
$ws.Range("D4").Value = '2024-11-22T14:18:00+00:00'
$ws.Range("F4").Value = 'KAR MARK kudampuli lehyam 1KG'
$ws.Range("I4").Value = 1230
$ws.Range("J4").Value = 131.78

$ws.Range("B5").Value = '406-5202604-5197163'
$ws.Range("C5").Value = '2024-11-22T10:44:10+00:00'
$ws.Range("D5").Value = '2024-11-22T15:03:38+00:00'
$ws.Range("J5").Value = 131.78
$ws.Range("L5").Value = 8.58

$ws.Range("B6").Value = '407-0749928-3148341'
$ws.Range("C6").Value = '2024-11-22T09:15:16+00:00'
$ws.Range("D6").Value = '2024-11-22T09:30:28+00:00'
$ws.Range("F6").Value = 'KAR MARK kudampuli lehyam 1KG'
$ws.Range("I6").Value = 1230
$ws.Range("J6").Value = 131.79
$ws.Range("K6").Value = 40
$ws.Range("L6").Value = 4.29

$ws.Range("B7").Value = '171-2860293-4573126'
$ws.Range("C7").Value = '2024-11-22T07:08:41+00:00'
$ws.Range("D7").Value = '2024-11-22T07:23:47+00:00'

$ws.Range("B8").Value = '402-7761907-9018739'
$ws.Range("C8").Value = '2024-11-22T21:07:55+00:00'
$ws.Range("D8").Value = '2024-11-22T21:21:15+00:00'
$ws.Range("F8").Value = 'KAR Mark Kudampuli Dried Lehyam (500 Gm)'
$ws.Range("I8").Value = 640
$ws.Range("J8").Value = 68.58

$ws.Range("B9").Value = '408-3619872-6041103'
$ws.Range("C9").Value = '2024-11-22T20:33:18+00:00'
$ws.Range("D9").Value = '2024-11-23T15:10:37+00:00'
$ws.Range("J9").Value = 68.58
$ws.Range("L9").Value = 8.58

$ws.Range("B10").Value = '407-5062014-7214742'
$ws.Range("C10").Value = '2024-11-22T17:51:42+00:00'
$ws.Range("D10").Value = '2024-11-22T18:07:04+00:00'
$ws.Range("F10").Value = 'KAR Mark Kudampuli Dried Lehyam (500 Gm)'
$ws.Range("I10").Value = 640
$ws.Range("J10").Value = 68.58

$ws.Range("B11").Value = '405-2890357-6694768'
$ws.Range("C11").Value = '2024-11-22T15:55:43+00:00'
$ws.Range("D11").Value = '2024-11-22T16:10:50+00:00'

$ws.Range("B12").Value = '408-2937479-4182750'
$ws.Range("C12").Value = '2024-11-22T13:05:03+00:00'
$ws.Range("D12").Value = '2024-11-22T13:20:13+00:00'
$ws.Range("J12").Value = 68.56999999999999
$ws.Range("L12").Value = 8.57

$ws.Range("B13").Value = '407-6478206-4772334'
$ws.Range("C13").Value = '2024-11-22T09:54:02+00:00'
$ws.Range("D13").Value = '2024-11-22T10:03:29+00:00'
$ws.Range("J13").Value = 68.58
$ws.Range("K13").Value = 80
$ws.Range("L13").Value = 8.58

$ws.Range("B14").Value = '171-7906369-9845952'
$ws.Range("C14").Value = '2024-11-22T09:32:17+00:00'
$ws.Range("D14").Value = '2024-11-22T09:41:33+00:00'
$ws.Range("F14").Value = 'KAR Mark Kudampuli Dried Lehyam (500 Gm)'
$ws.Range("I14").Value = 640
$ws.Range("J14").Value = 68.58
$ws.Range("K14").Value = 80
$ws.Range("L14").Value = 8.58

$ws.Range("B15").Value = '407-0749928-3148341'
$ws.Range("C15").Value = '2024-11-22T09:15:16+00:00'
$ws.Range("D15").Value = '2024-11-22T09:30:28+00:00'
$ws.Range("J15").Value = 68.56999999999999
$ws.Range("K15").Value = 40
$ws.Range("L15").Value = 4.29

$ws.Range("B16").Value = '405-4471990-3110724'
$ws.Range("C16").Value = '2024-11-22T08:29:42+00:00'
$ws.Range("D16").Value = '2024-11-22T08:44:25+00:00'

$ws.Range("B17").Value = '405-5684197-8263504'
$ws.Range("C17").Value = '2024-11-22T08:15:07+00:00'
$ws.Range("D17").Value = '2024-11-22T08:30:16+00:00'
$ws.Range("F17").Value = 'KAR Mark Kudampuli Dried Lehyam (500 Gm)'
$ws.Range("I17").Value = 640
$ws.Range("J17").Value = 68.58

$ws.Range("B18").Value = '405-5142500-8749156'
$ws.Range("C18").Value = '2024-11-22T07:26:12+00:00'
$ws.Range("D18").Value = '2024-11-22T07:43:38+00:00'
$ws.Range("F18").Value = 'KAR Mark Kudampuli Dried Lehyam For Weight Loss And Reduce Belly Fat(750 Gm)'
$ws.Range("I18").Value = 950
$ws.Range("J18").Value = 101.78
